# Fix Training Data Issue (#48)
# The "Date" column (BF) for every data row (2-31) on Sheet1 was off by one
# day ("5-8-2007-08") because of how the NBA stats export formatted the
# game date. Correct it to the proper ISO-style date text "2008-05-08".
#
# Column BF is the 58th column.  A leading apostrophe is used so Excel
# stores the corrected value as literal text (matching the original cell's
# text type) instead of auto-converting the "2008-05-08" pattern into a
# real date serial number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateColumn = 58  # column BF
$firstDataRow = 2
$lastDataRow = 31
$correctedDate = "'2008-05-08"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $dateColumn).Value = $correctedDate
}
